# The edit rotates the "species observation" data cyclically among rows 3-6:
#   new row3 <- old row5
#   new row4 <- old row6
#   new row5 <- old row4
#   new row6 <- old row3
# Only the data columns (A,B,E,F,G,H,Q,R) and the sparse/empty cells in the
# I:N / AF block actually move; every other column (C,D,P,S,T,U,V,W,Y,Z,AA,AB,
# AD,AE,AG,AT,AW,AX,AY) is identical across these four rows and is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- capture the "before" state of the rotating value columns ----
$cols = @("A","B","E","F","G","H","Q","R")

$vals = @{}
foreach ($r in 3,4,5,6) {
    $vals[$r] = @{}
    foreach ($c in $cols) {
        $vals[$r][$c] = $ws.Range("$c$r").Value()
    }
}

# permutation: new row $r gets the old content of row $src[$r]
$src = @{ 3 = 5; 4 = 6; 5 = 4; 6 = 3 }

foreach ($r in 3,4,5,6) {
    $s = $src[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $vals[$s][$c]
    }
}

# ---- row 3: M3 ("färska spår") goes away, J3 becomes a new empty cell, AF3 becomes a new empty cell ----
$ws.Range("M3").Value = $null
$ws.Range("I3").Copy($ws.Range("J3"))
$ws.Range("I3").Copy($ws.Range("AF3"))

# ---- row 4: L4 becomes a new empty cell ----
$ws.Range("I4").Copy($ws.Range("L4"))

# ---- row 5: L5 is removed ----
$ws.Range("L5").Value = $null

# ---- row 6: J6 and AF6 are removed, M6 becomes present with "färska spår" ----
$ws.Range("J6").Value = $null
$ws.Range("AF6").Value = $null
$ws.Range("M6").Value = "färska spår"
